# Apply updated crypto price/volume figures (GitHub Actions data refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.501.31"
$ws.Range("E2").Value = "  +0.33%  "
$ws.Range("D3").Value = "1.913.65"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.82"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4832"
$ws.Range("E7").Value = "  +3.33%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2892"
$ws.Range("E8").Value = "  +1.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06702"
$ws.Range("E9").Value = "  -1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "109.66"
$ws.Range("E10").Value = "  +2.69%  "
$ws.Range("E11").Value = "  +6.25%  "
$ws.Range("D12").Value = "1.908.75"
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07557"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.266"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6683"
$ws.Range("E15").Value = "  +2.52%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "275.23"
$ws.Range("E16").Value = "  -3.50%  "
$ws.Range("D17").Value = "30.491.90"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9995"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007530"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.85"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").Value = "2.163.23"
$ws.Range("E21").Value = "  +0.31%  "
$ws.Range("E22").Value = "  +5.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.03%  "
$ws.Range("E24").Value = "  +4.40%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.437"
$ws.Range("E25").Value = "  +2.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.52"
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.19"
$ws.Range("E27").Value = "  -5.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.122"
$ws.Range("E28").Value = "  +4.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1051"
$ws.Range("E29").Value = "  -1.66%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.402"
$ws.Range("E30").Value = "  +2.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.141"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.052"
$ws.Range("E32").Value = "  +3.06%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04993"
$ws.Range("E33").Value = "  -0.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7285"
$ws.Range("E34").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9990"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.730"
$ws.Range("E38").Value = "  +1.13%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.669"
$ws.Range("E39").Value = "  -0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "110.85"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  -0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4421"
$ws.Range("E42").Value = "  +5.69%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8660"
$ws.Range("E43").Value = "  -0.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.848"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "67.89"
$ws.Range("E46").Value = "  +1.10%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.424"
$ws.Range("E47").Value = "  +4.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.284"
$ws.Range("E48").Value = "  +1.16%  "
$ws.Range("E49").Value = "  +3.35%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "47.60"
$ws.Range("E50").Value = "  -9.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.470"
$ws.Range("E51").Value = "  +7.62%  "
